$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the "Full Offense Details" column (old column L)
# ---------------------------------------------------------------------------
$ws.Range("L1").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2. Insert two new columns before "Jail Sentence" (now column N) for the
#    new "Payments" / "Combined Fees/Fines/Costs" fields.
# ---------------------------------------------------------------------------
$ws.Range("N1:O1").EntireColumn.Insert()
$ws.Range("N1").Value = "Payments"
$ws.Range("O1").Value = "Combined Fees/Fines/Costs"

# ---------------------------------------------------------------------------
# 3. Rename "Jail Sentence" (now column P) to "Jail Sentence(s)" and append
#    the new trailing columns.
# ---------------------------------------------------------------------------
$ws.Range("P1").Value = "Jail Sentence(s)"
$ws.Range("Q1").Value = "Sentence Years"
$ws.Range("R1").Value = "Raw Case Summary"
$ws.Range("S1").Value = "Raw Offense Information"
$ws.Range("T1").Value = "Raw Payments Made to the Court"
$ws.Range("U1").Value = "Raw Register of Actions"

# ---------------------------------------------------------------------------
# 4. Column widths
# ---------------------------------------------------------------------------
$ws.Range("A1:B1").ColumnWidth = 16.17
$ws.Range("C1").ColumnWidth = 13.12
$ws.Range("D1").ColumnWidth = 21.04
$ws.Range("E1").ColumnWidth = 9.36
$ws.Range("F1").ColumnWidth = 16.73
$ws.Range("G1").ColumnWidth = 10.19
$ws.Range("H1").ColumnWidth = 10.75
$ws.Range("I1").ColumnWidth = 11.17
$ws.Range("J1").ColumnWidth = 17.94
$ws.Range("K1").ColumnWidth = 17.8
$ws.Range("L1").ColumnWidth = 17.52
$ws.Range("N1").ColumnWidth = 12.15
$ws.Range("O1").ColumnWidth = 25.22
$ws.Range("P1").ColumnWidth = 18.35
$ws.Range("Q1").ColumnWidth = 13.95
$ws.Range("R1").ColumnWidth = 40.75
$ws.Range("S1").ColumnWidth = 40.05
$ws.Range("T1").ColumnWidth = 41.03
$ws.Range("U1").ColumnWidth = 57.72

# ---------------------------------------------------------------------------
# 5. Header row styling: bold font + thin bottom border across A1:U1
# ---------------------------------------------------------------------------
$headerMain = $ws.Range("A1:Q1")
$headerMain.Font.Bold = $true
$headerMain.Borders.Item(9).LineStyle = 1

$headerRaw = $ws.Range("R1:U1")
$headerRaw.Font.Bold = $true
$headerRaw.Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------------------
# 6. View bits (top-left cell / selection) to match the authored state
# ---------------------------------------------------------------------------
$null = $ws.Range("K2").Select()
$excel.ActiveWindow.ScrollColumn = 8

Write-Host "done"
